# Update F-column ("想去人数" / interest counts) values to match refreshed
# scrape data, per commit "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$updates = @{
    2 = 181
    3 = 3123
    4 = 843
    5 = 614
    6 = 610
    7 = 489
    8 = 666
    10 = 637
    12 = 2182
    13 = 1278
    14 = 762
    15 = 59
    17 = 2713
    19 = 38
    21 = 541
    23 = 757
    24 = 44
    25 = 34
    27 = 41
    28 = 574
    29 = 573
    32 = 927
    33 = 4764
    34 = 307
    35 = 70
    36 = 29
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item("演出")
$updates = @{
    7 = 377
    14 = 280
    16 = 63
    21 = 13
    24 = 19
    27 = 384
    37 = 639
    38 = 639
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item("本地生活")
$updates = @{
    4 = 1485
    5 = 591
    6 = 325
    7 = 304
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    3 = 1485
    4 = 591
    5 = 181
    6 = 325
    7 = 3123
    8 = 843
    9 = 614
    10 = 610
    11 = 489
    12 = 666
    14 = 637
    19 = 2182
    20 = 1278
    21 = 762
    22 = 59
    24 = 2713
    26 = 38
    27 = 541
    28 = 63
    29 = 304
    32 = 757
    33 = 757
    34 = 44
    35 = 13
    36 = 34
    37 = 41
    39 = 19
    40 = 574
    43 = 927
    44 = 4764
    46 = 307
    48 = 70
    50 = 639
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
